$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Лист1")
$ws1.Name = "2D KERNEL VALUES"
$ws2 = $wb.Worksheets.Item("Лист2")
$ws2.Name = "2D KERNEL SURFACE"
